# ---------------------------------------------------------------------------
# Rename the two existing sheets and add a new "Datos" sheet, refresh the
# "saldo" (H) column on the example data sheet, fix a mis-typed category on
# one row, and update the saved selections / active tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Rename sheets: Hoja2 -> Intrucciones, Hoja1 -> Ejemplo -----------------
$wsInstrucciones = $wb.Worksheets.Item(1)
$wsInstrucciones.Name = "Intrucciones"

$wsEjemplo = $wb.Worksheets.Item(2)
$wsEjemplo.Name = "Ejemplo"

# --- Fix a mis-categorised product row: "BICLETAS" (typo) -> "CASCOS" ------
$wsEjemplo.Range("C35").Value = "CASCOS"

# --- Fill in the new "saldo" (column H) values for data rows 2-53 ----------
$saldoData = @(
    @(2,40),   @(3,51),   @(4,34),   @(5,45),   @(6,40),
    @(7,58),   @(8,55),   @(9,24),   @(10,48),  @(11,20),
    @(12,50),  @(13,32),  @(14,34),  @(15,33),  @(16,49),
    @(17,60),  @(18,23),  @(19,46),  @(20,56),  @(21,46),
    @(22,32),  @(23,44),  @(24,31),  @(25,25),  @(26,53),
    @(27,45),  @(28,31),  @(29,39),  @(30,32),  @(31,45),
    @(32,22),  @(33,42),  @(34,20),  @(35,18),  @(36,29),
    @(37,27),  @(38,17),  @(39,26),  @(40,22),  @(41,27),
    @(42,41),  @(43,44),  @(44,30),  @(45,42),  @(46,25),
    @(47,47),  @(48,45),  @(49,22),  @(50,17),  @(51,14),
    @(52,27),  @(53,14)
)
foreach ($pair in $saldoData) {
    $r = $pair[0]
    $v = $pair[1]
    $wsEjemplo.Cells.Item($r, 8).Value = $v
}

# --- Restore the per-sheet selection state ----------------------------------
[void]$wsInstrucciones.Range("C11").Select()
[void]$wsEjemplo.Range("A1:H1").Select()

# --- Add the new "Datos" sheet at the end, seeded with the same header -----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDatos = $wb.Worksheets.Add($null, $lastSheet)
$wsDatos.Name = "Datos"
[void]$wsEjemplo.Range("A1:H1").Copy($wsDatos.Range("A1"))
[void]$wsDatos.Range("F9").Select()

# --- "Datos" is the active sheet/tab when the file is saved -----------------
[void]$wsDatos.Activate()
